# Applies the "gh-pages output regenerated" update:
#  - Sheet "展览"   : remove the duplicated Gumi-doujin-show row (old row 17),
#                     then zero out every "想去人数" (F) value.
#  - Sheet "演出"   : zero out every "想去人数" (F) value (no row shift here).
#  - Sheet "本地生活": header only, untouched.
#  - Sheet "全部类型": remove the duplicated Gumi-doujin-show row (old row 18),
#                     then zero out every "想去人数" (F) value.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -----------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
# Row 17 duplicates row 16 ("合肥·首届Gumi同人展"); delete it so row 18's
# "巢湖·喵喵漫游戏动漫展" (and everything after) shifts up by one.
$ws1.Rows.Item(17).Delete()
# Zero every "想去人数" value in the data rows (now 2..27 after the delete).
$ws1.Range("F2:F27").Value = 0

# --- Sheet "演出" -----------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2:F3").Value = 0

# --- Sheet "本地生活" --------------------------------------------------
# Header-only sheet; no data rows to update.

# --- Sheet "全部类型" --------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
# Row 18 duplicates row 17 ("合肥·首届Gumi同人展"); delete it so row 19's
# "巢湖·喵喵漫游戏动漫展" (and everything after) shifts up by one.
$ws4.Rows.Item(18).Delete()
# Zero every "想去人数" value in the data rows (now 2..29 after the delete).
$ws4.Range("F2:F29").Value = 0
